$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stats")

# Row 2 (knight): branch id 1 -> 0, jumppower 7 -> 15
$ws.Range("A2").Value = 0
$ws.Range("D2").Value = 15

# Row 3: branch id 2 -> 1, name cat -> rogue
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "rogue"

# Row 4: branch id 3 -> 2, name mage -> dog
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "dog"

# Update selection to B4
$ws.Range("B4").Select()
